$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.274482131004333
$ws.Range("B1").Value = 2.480945348739624
$ws.Range("C1").Value = 3.480561971664429
$ws.Range("D1").Value = 3.09826922416687
$ws.Range("E1").Value = 1.067774653434753
